# Actualización automática 2025-10-20 17:30:09
$wb = $excel.ActiveWorkbook

# Sheet: VENTAS POR GRUPO
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M22").Value = 9981.23

# Sheet: VENTA MENSUAL
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F22").Value = 11593.25
$ws2.Range("F26").Value = 23462.9

# Sheet: CUMPLIMIENTO MENSUAL
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D12").Value = 22455.08
$ws3.Range("E12").Value = 5499.899999999998
$ws3.Range("F12").Value = 0.8032586680441196
$ws3.Range("D14").Value = 22292.69
$ws3.Range("E14").Value = 19910.69110009468
$ws3.Range("F14").Value = 0.5282204747322955
